# Daily attendance processing - 2026-01-19 08:05:12
# Normalize the "Recorded By" column so that the System entry is listed
# after the user's email instead of before it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Session Analysis Results")

$ws.Cells.Replace(
    "System, dnasr281@gmail.com",
    "dnasr281@gmail.com, System",
    1,
    1,
    $false,
    $false,
    $false,
    $false
)
